$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 71") {
        $target = $sh
        break
    }
}

if ($target -ne $null) {
    $target.TextFrame.TextRange.Text = "Azure`rSynapse`rAnalytics"
}
